$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-7 with new recorded values
$ws.Range("B2").Value = -0.003550
$ws.Range("C2").Value = 2.239148
$ws.Range("D2").Value = 0.000000

$ws.Range("B3").Value = -0.007835
$ws.Range("C3").Value = 2.288220
$ws.Range("D3").Value = 0.049072

$ws.Range("B4").Value = -0.012120
$ws.Range("C4").Value = 2.288220
$ws.Range("D4").Value = 0.049072

$ws.Range("B5").Value = -0.016405
$ws.Range("C5").Value = 2.239148
$ws.Range("D5").Value = 0.000000

$ws.Range("B6").Value = -0.020690
$ws.Range("C6").Value = 2.484510
$ws.Range("D6").Value = 0.245362

$ws.Range("B7").Value = -0.024975
$ws.Range("C7").Value = 2.828018
$ws.Range("D7").Value = 0.588870

# Add new row 8 with data, copying style from row 7
$ws.Range("A7:D7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = 6.000000
$ws.Range("B8").Value = -0.029260
$ws.Range("C8").Value = 3.171525
$ws.Range("D8").Value = 0.932377
